# Updates the cryptos list: refresh Price (col D) and Volume(1h) (col E) values,
# and shift the coin entries from row 36 onward by one position (RenzoRestakedETH
# drops off the list, and Bittensor is appended at the end). Price values are
# written with a leading apostrophe so Excel keeps them as text (matching the
# workbook's original formatting) instead of silently converting numeric-
# looking strings like "1.00" or "604.96" into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.792.22"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "'3.820.31"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'604.96"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").Value = "'166.96"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").Value = "'0.452"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "'4.458.15"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'3.812.21"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "'18.49"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "'67.812.36"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "'7.10"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "'463.06"
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("D21").Value = "'9.92"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").Value = "'0.702"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").Value = "'83.43"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'12.07"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'2.12"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "'10.10"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "'3.966.27"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "'29.67"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "'9.09"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.100"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.37"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.138"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'5.81"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'48.13"
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.301"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'28.33"
$ws.Range("E45").Value = "  +10.90%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "'43.21"
$ws.Range("E46").Value = "  -5.07%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.40"
$ws.Range("E47").Value = "  +11.90%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'8.35"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'148.20"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.84"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'389.28"
$ws.Range("E51").Value = "  -0.14%  "
